$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.344.38'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.841.86'
$ws.Range('E3').Value = '  -0.44%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.9999'
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '240.50'
$ws.Range('E5').Value = '  -0.09%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.6273'
$ws.Range('E6').Value = '  -0.18%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.001'
$ws.Range('E7').Value = '  +0.02%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.07480'
$ws.Range('E8').Value = '  -2.78%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.2892'
$ws.Range('E9').Value = '  -1.10%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '24.34'
$ws.Range('E10').Value = '  -2.56%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07715'
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').Value = '1.841.90'
$ws.Range('E12').Value = '  -0.98%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '4.986'
$ws.Range('E13').Value = '  -1.09%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.6778'
$ws.Range('E14').Value = '  -0.71%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.00001029'
$ws.Range('E15').Value = '  -4.92%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '82.08'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = '2.103.89'
$ws.Range('E17').Value = '  -0.94%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '6.121'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = '29.363.72'
$ws.Range('E19').Value = '  -0.32%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '228.84'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  -1.16%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '1.001'
$ws.Range('E22').Value = '  +0.00%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '7.370'
$ws.Range('E23').Value = '  -1.37%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '1.002'
$ws.Range('E24').Value = '  +0.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '158.85'
$ws.Range('E25').Value = '  +0.64%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.1381'
$ws.Range('E26').Value = '  +0.03%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.381'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -1.23%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.399'
$ws.Range('E29').Value = '  +2.60%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.475'
$ws.Range('E30').Value = '  +0.86%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.05680'
$ws.Range('E31').Value = '  +0.98%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.094'
$ws.Range('E32').Value = '  -0.82%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.042'
$ws.Range('E33').Value = '  -0.20%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.818'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('E35').Value = '  -1.71%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.6960'
$ws.Range('E36').Value = '  -1.34%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.586'
$ws.Range('E37').Value = '  -0.44%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.838'
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('D39').Value = '1.251.59'
$ws.Range('E39').Value = '  +2.04%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.01813'
$ws.Range('E40').Value = '  +1.18%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '6.515'
$ws.Range('E41').Value = '  +0.90%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.9033'
$ws.Range('E42').Value = '  -0.01%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.9999'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').Value = '2.005.89'
$ws.Range('E44').Value = '  -1.18%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '101.10'
$ws.Range('E45').Value = '  -0.90%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '65.76'
$ws.Range('E46').Value = '  -0.60%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '7.069'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.00000000117'
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.1162'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '8.953'
$ws.Range('E50').Value = '  -1.01%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.3936'
$ws.Range('E51').Value = '  -2.16%  '
